$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 2.179726059491315
$ws.Cells.Item(2, 3).Value = 0.1228732742893044
$ws.Cells.Item(2, 4).Value = 0.2304411993544448
$ws.Cells.Item(2, 5).Value = 0.05807592456839572
$ws.Cells.Item(2, 6).Value = 3.901996374871942
$ws.Cells.Item(2, 12).Value = 0.2541781385875055
$ws.Cells.Item(3, 2).Value = 2.089191184859203
$ws.Cells.Item(3, 3).Value = 0.1063440053391673
$ws.Cells.Item(3, 4).Value = 0.220427376401247
$ws.Cells.Item(3, 5).Value = 0.05772313933465512
$ws.Cells.Item(3, 6).Value = 3.702702403507089
$ws.Cells.Item(3, 12).Value = 0.2461069628529629
$ws.Cells.Item(4, 2).Value = 2.035691974967222
$ws.Cells.Item(4, 3).Value = 0.0962401253910059
$ws.Cells.Item(4, 4).Value = 0.2142908285513698
$ws.Cells.Item(4, 5).Value = 0.05751285478061874
$ws.Cells.Item(4, 6).Value = 3.581514854234285
$ws.Cells.Item(4, 12).Value = 0.2413381715008285
$ws.Cells.Item(5, 2).Value = 2.014412456097944
$ws.Cells.Item(5, 3).Value = 0.0921332781626063
$ws.Cells.Item(5, 4).Value = 0.2117925400974059
$ws.Cells.Item(5, 5).Value = 0.05742872881557304
$ws.Cells.Item(5, 6).Value = 3.532417022318725
$ws.Cells.Item(5, 12).Value = 0.239441537443227
$ws.Cells.Item(6, 2).Value = 2.010910423114012
$ws.Cells.Item(6, 3).Value = 0.09145195467456801
$ws.Cells.Item(6, 4).Value = 0.211377828750912
$ws.Cells.Item(6, 5).Value = 0.05741485368395338
$ws.Cells.Item(6, 6).Value = 3.524281425495758
$ws.Cells.Item(6, 12).Value = 0.2391294133544761
$ws.Cells.Item(7, 2).Value = 2.035402883674635
$ws.Cells.Item(7, 3).Value = 0.09618469726524381
$ws.Cells.Item(7, 4).Value = 0.2142571268249185
$ws.Cells.Item(7, 5).Value = 0.05751171391678689
$ws.Cells.Item(7, 6).Value = 3.580851553921462
$ws.Cells.Item(7, 12).Value = 0.241312404213744
$ws.Cells.Item(8, 2).Value = 2.148073715398198
$ws.Cells.Item(8, 3).Value = 0.1171641256604516
$ws.Cells.Item(8, 4).Value = 0.2269855772361637
$ws.Cells.Item(8, 5).Value = 0.05795295624357522
$ws.Cells.Item(8, 6).Value = 3.833029101149691
$ws.Cells.Item(8, 12).Value = 0.2513561983885779
$ws.Cells.Item(9, 2).Value = 2.385782493052375
$ws.Cells.Item(9, 3).Value = 0.158701129582397
$ws.Cells.Item(9, 4).Value = 0.2520696190262015
$ws.Cells.Item(9, 5).Value = 0.05886958666623165
$ws.Cells.Item(9, 6).Value = 4.337368848030678
$ws.Cells.Item(9, 12).Value = 0.2725515286154945
$ws.Cells.Item(10, 2).Value = 2.570925862391903
$ws.Cells.Item(10, 3).Value = 0.1895164179622668
$ws.Cells.Item(10, 4).Value = 0.2706153454228968
$ws.Cells.Item(10, 5).Value = 0.05957602564258835
$ws.Cells.Item(10, 6).Value = 4.714570110744347
$ws.Cells.Item(10, 12).Value = 0.2890626982863296
$ws.Cells.Item(11, 2).Value = 2.657496254550779
$ws.Cells.Item(11, 3).Value = 0.2036124836000965
$ws.Cells.Item(11, 4).Value = 0.2790866267354488
$ws.Cells.Item(11, 5).Value = 0.05990494031451554
$ws.Cells.Item(11, 6).Value = 4.887765828894601
$ws.Cells.Item(11, 12).Value = 0.29678359886816
$ws.Cells.Item(12, 2).Value = 2.69062076682377
$ws.Cells.Item(12, 3).Value = 0.2089625199811849
$ws.Cells.Item(12, 4).Value = 0.2823001902308988
$ws.Cells.Item(12, 5).Value = 0.06003060721674203
$ws.Cells.Item(12, 6).Value = 4.953593037336759
$ws.Cells.Item(12, 12).Value = 0.2997379184984936
$ws.Cells.Item(13, 2).Value = 2.68347151721224
$ws.Cells.Item(13, 3).Value = 0.2078097377757047
$ws.Cells.Item(13, 4).Value = 0.2816078286330992
$ws.Cells.Item(13, 5).Value = 0.06000349260095206
$ws.Cells.Item(13, 6).Value = 4.939405041871169
$ws.Cells.Item(13, 12).Value = 0.2991002862563619
$ws.Cells.Item(14, 2).Value = 2.660214543484358
$ws.Cells.Item(14, 3).Value = 0.2040523850908187
$ws.Cells.Item(14, 4).Value = 0.2793508909241211
$ws.Cells.Item(14, 5).Value = 0.05991525651617025
$ws.Cells.Item(14, 6).Value = 4.893176564512942
$ws.Cells.Item(14, 12).Value = 0.2970260372742501
$ws.Cells.Item(15, 2).Value = 2.646013676394716
$ws.Cells.Item(15, 3).Value = 0.2017525102867523
$ws.Cells.Item(15, 4).Value = 0.2779692104936089
$ws.Cells.Item(15, 5).Value = 0.0598613552913303
$ws.Cells.Item(15, 6).Value = 4.864892098538348
$ws.Cells.Item(15, 12).Value = 0.295759494666072
$ws.Cells.Item(16, 2).Value = 2.565315914551661
$ws.Cells.Item(16, 3).Value = 0.1885968570958596
$ws.Cells.Item(16, 4).Value = 0.2700624840178136
$ws.Cells.Item(16, 5).Value = 0.05955468465093716
$ws.Cells.Item(16, 6).Value = 4.703284591627494
$ws.Cells.Item(16, 12).Value = 0.2885623767262615
$ws.Cells.Item(17, 2).Value = 2.51641497859498
$ws.Cells.Item(17, 3).Value = 0.1805469389481686
$ws.Cells.Item(17, 4).Value = 0.2652213400374137
$ws.Cells.Item(17, 5).Value = 0.05936850623124101
$ws.Cells.Item(17, 6).Value = 4.60456242849321
$ws.Cells.Item(17, 12).Value = 0.2842012165261849
$ws.Cells.Item(18, 2).Value = 2.488509151719484
$ws.Cells.Item(18, 3).Value = 0.1759241264217053
$ws.Cells.Item(18, 4).Value = 0.2624400706334598
$ws.Cells.Item(18, 5).Value = 0.0592621304535399
$ws.Cells.Item(18, 6).Value = 4.547930452495763
$ws.Cells.Item(18, 12).Value = 0.2817125234062985
$ws.Cells.Item(19, 2).Value = 2.479098486112093
$ws.Cells.Item(19, 3).Value = 0.1743601448848153
$ws.Cells.Item(19, 4).Value = 0.2614989146065057
$ws.Cells.Item(19, 5).Value = 0.05922623438877395
$ws.Cells.Item(19, 6).Value = 4.528781339271575
$ws.Cells.Item(19, 12).Value = 0.2808732714517248
$ws.Cells.Item(20, 2).Value = 2.521597692578609
$ws.Cells.Item(20, 3).Value = 0.181403105509844
$ws.Cells.Item(20, 4).Value = 0.2657363490518492
$ws.Cells.Item(20, 5).Value = 0.05938825164811945
$ws.Cells.Item(20, 6).Value = 4.615055926583352
$ws.Cells.Item(20, 12).Value = 0.2846634245373139
$ws.Cells.Item(21, 2).Value = 2.66703635642574
$ws.Cells.Item(21, 3).Value = 0.2051556729606432
$ws.Cells.Item(21, 4).Value = 0.2800136488241947
$ws.Cells.Item(21, 5).Value = 0.05994114313185506
$ws.Cells.Item(21, 6).Value = 4.906748339658407
$ws.Cells.Item(21, 12).Value = 0.2976344612404915
$ws.Cells.Item(22, 2).Value = 2.764085669484757
$ws.Cells.Item(22, 3).Value = 0.2207506932515173
$ws.Cells.Item(22, 4).Value = 0.2893780780115662
$ws.Cells.Item(22, 5).Value = 0.06030899662811606
$ws.Cells.Item(22, 6).Value = 5.098800512015487
$ws.Cells.Item(22, 12).Value = 0.306290233872204
$ws.Cells.Item(23, 2).Value = 2.712104371091243
$ws.Cells.Item(23, 3).Value = 0.2124204924413675
$ws.Cells.Item(23, 4).Value = 0.2843768244978833
$ws.Cells.Item(23, 5).Value = 0.06011206150081705
$ws.Cells.Item(23, 6).Value = 4.99616557335878
$ws.Cells.Item(23, 12).Value = 0.301654021518786
$ws.Cells.Item(24, 2).Value = 2.519253940143017
$ws.Cells.Item(24, 3).Value = 0.1810160165042589
$ws.Cells.Item(24, 4).Value = 0.2655035073959198
$ws.Cells.Item(24, 5).Value = 0.05937932268940038
$ws.Cells.Item(24, 6).Value = 4.610311428477075
$ws.Cells.Item(24, 12).Value = 0.2844544024032416
$ws.Cells.Item(25, 2).Value = 2.319653025489401
$ws.Cells.Item(25, 3).Value = 0.1474160792265025
$ws.Cells.Item(25, 4).Value = 0.2452662655368698
$ws.Cells.Item(25, 5).Value = 0.05861594346739274
$ws.Cells.Item(25, 6).Value = 4.199810963210666
$ws.Cells.Item(25, 12).Value = 0.2666545708321735
